# Part 8 completed. Project runs. Drop menu for items and functionality added.
#
# The paragraph that used to introduce the creation of "menus.py" is
# rewritten into a short sentence that instead introduces the drop-menu
# reuse of the existing inventory-menu function. That paragraph currently
# holds three runs:
#   1) "Before we get into ... creating a new file, called "   (Segoe UI)
#   2) "menus.py"                                              (HTMLCode / Consolas)
#   3) ", where we'll store ... in that file:"                 (Segoe UI)
# After the edit it must hold a single run (reusing run 1's formatting)
# whose text is the new sentence; runs 2 and 3 disappear entirely.

$d = $word.ActiveDocument
$rsq = [char]0x2019

$newText = "Now for displaying the drop menu. It" + $rsq + `
    "s really not different from the inventory menu, so we can use the " + `
    "same function, and send a different title to it."

# Locate the start of run 1 via a stable anchor phrase near its beginning.
$anchorRange = $d.Content.Duplicate
$gotAnchor = $anchorRange.Find.Execute("Before we get into how to use items")
if (-not $gotAnchor) {
    throw "Could not find the paragraph to edit."
}
$startPos = $anchorRange.Start

# Locate the start of run 2 ("menus.py"), i.e. the end of run 1.
$codeRange = $d.Content.Duplicate
$gotCode = $codeRange.Find.Execute("menus.py")
if (-not $gotCode) {
    throw "Could not find the 'menus.py' run."
}
$codeStart = $codeRange.Start

# Step 1: overwrite run 1's text with the new sentence (formatting of run 1
# is preserved since we are only replacing the text of its own range).
$firstRun = $d.Range($startPos, $codeStart)
$firstRun.Text = $newText

# Step 2: remove the old "menus.py" run together with the old trailing run
# ("... Put the following code in that file:"), which now immediately
# follows the freshly-inserted sentence.
$tailRange = $d.Content.Duplicate
$gotTail = $tailRange.Find.Execute("Put the following code in that file:")
if (-not $gotTail) {
    throw "Could not find the trailing run to remove."
}
$tailEnd = $tailRange.End

$newTextEnd = $startPos + $newText.Length
$toRemove = $d.Range($newTextEnd, $tailEnd)
$toRemove.Delete()
